$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 362
$ws.Range("I2").Value = 1066
$ws.Range("J2").Value = 4484
$ws.Range("K2").Value = 21
$ws.Range("L2").Value = 1232
$ws.Range("M2").Value = 59
$ws.Range("N2").Value = 716
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 26
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 58
$ws.Range("S2").Value = 462
$ws.Range("T2").Value = 788
$ws.Range("U2").Value = 53
$ws.Range("V2").Value = 6969
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 7066
$ws.Range("Y2").Value = 16
$ws.Range("Z2").Value = 109
$ws.Range("AA2").Value = 54

$wb.Save()
